# Weekly update: insert two new rows of "Acelga" price data at the top of
# the existing date-ordered block (rows 413-414), pushing the previously
# existing rows 413-434 down to 415-436 (row 435 -> 437 unchanged in content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 413; this shifts rows
# 413..435 down to 415..437, preserving all of their existing values/styles.
$ws.Range("A413:A414").EntireRow.Insert()

# New row 413: "Primera" quality entry for 2023-04-25 (serial 45041)
$ws.Range("A413").Value = 7
$ws.Range("B413").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C413").Value = "Ñuble"
$ws.Range("D413").Value = 45041
$ws.Range("E413").Value = 16
$ws.Range("F413").Value = 100112009
$ws.Range("G413").Value = "Acelga"
$ws.Range("H413").Value = "Sin especificar"
$ws.Range("I413").Value = "Primera"
$ws.Range("J413").Value = 300
$ws.Range("K413").Value = 700
$ws.Range("L413").Value = 700
$ws.Range("M413").Value = 700
$ws.Range("N413").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O413").Value = "Provincia de Diguillín"
$ws.Range("P413").Value = 700
$ws.Range("Q413").Value = 1
$ws.Range("R413").Value = "Hortaliza"

# New row 414: "Segunda" quality entry for the same date 2023-04-25 (serial 45041)
$ws.Range("A414").Value = 7
$ws.Range("B414").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C414").Value = "Ñuble"
$ws.Range("D414").Value = 45041
$ws.Range("E414").Value = 16
$ws.Range("F414").Value = 100112009
$ws.Range("G414").Value = "Acelga"
$ws.Range("H414").Value = "Sin especificar"
$ws.Range("I414").Value = "Segunda"
$ws.Range("J414").Value = 300
$ws.Range("K414").Value = 500
$ws.Range("L414").Value = 500
$ws.Range("M414").Value = 500
$ws.Range("N414").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O414").Value = "Provincia de Diguillín"
$ws.Range("P414").Value = 500
$ws.Range("Q414").Value = 1
$ws.Range("R414").Value = "Hortaliza"
